# Update DateBase/orders/International Ever Green_2025-10-29.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# New order rows appended after existing data (rows 42-51)
$rows = @(
    @{ Row = 42; A = $null; C = "480_蝴蝶洋牡丹红_butterfly  Ranunculus_undefined_1bunch"; F = "5" },
    @{ Row = 43; A = "4";   C = "663_大丽花 乌梅子酱_undefined_undefined_5stems";            F = "5" },
    @{ Row = 44; A = $null; C = "649_洋牡丹樱花粉_undefined_undefined_1bunch";               F = "5" },
    @{ Row = 45; A = "1";   C = "721_银扇干花_undefined_undefined_1bunch";                   F = "10" },
    @{ Row = 46; A = $null; C = "664_大丽花 琳达宝贝_undefined_undefined_5stems";            F = "10" },
    @{ Row = 47; A = "2";   C = "177_国王日_Kings Day_Rosa rugosa Thunb._20stems";           F = "10" },
    @{ Row = 48; A = $null; C = "118_绣球老绿_Hydrangea Garden Lace_Hydrangea L._1stem";     F = "20" },
    @{ Row = 49; A = $null; C = "592_进口春兰叶_undefined_undefined_1bunch";                 F = "10" },
    @{ Row = 50; A = $null; C = "540_糖棉_gomphocarpus fruticosus_undefined_1bunch";         F = $null },
    @{ Row = 51; A = $null; C = "816_山里红_undefined_undefined_1bunch";                     F = $null }
)

foreach ($r in $rows) {
    if ($r.A -ne $null) {
        $cell = $ws.Cells.Item($r.Row, 1)
        $cell.NumberFormat = "@"
        $cell.Value = $r.A
    }
    # FlowerName (column C) values are never purely numeric, so they are
    # already stored as text without needing an explicit text format.
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    if ($r.F -ne $null) {
        $cell = $ws.Cells.Item($r.Row, 6)
        $cell.NumberFormat = "@"
        $cell.Value = $r.F
    }
}

# Update Summary sheet G2 tracking code (string appended)
$ws2 = $wb.Worksheets.Item("Summary")
$cell = $ws2.Cells.Item(2, 7)
$cell.NumberFormat = "@"
$cell.Value = "02424158281012115205540502050501010156121014106165105105730101035255555101010201000"
